# Update the division problems in the table to match the new generated output.
$d = $word.ActiveDocument

$replacements = @(
    @("691÷9=", "518÷5="),
    @("629÷5=", "871÷3="),
    @("552÷8=", "816÷4="),
    @("979÷3=", "750÷7="),
    @("870÷3=", "648÷4="),
    @("947÷2=", "356÷8="),
    @("167÷3=", "459÷6="),
    @("842÷2=", "704÷7="),
    @("674÷5=", "154÷5="),
    @("961÷6=", "755÷8="),
    @("885÷4=", "887÷5="),
    @("402÷5=", "506÷4="),
    @("253÷4=", "491÷9="),
    @("670÷7=", "721÷8="),
    @("359÷9=", "638÷5="),
    @("737÷9=", "668÷4="),
    @("691÷5=", "703÷3="),
    @("297÷5=", "505÷7="),
    @("554÷9=", "466÷5="),
    @("683÷3=", "262÷3="),
    @("604÷9=", "637÷6="),
    @("735÷7=", "562÷2="),
    @("826÷6=", "482÷3="),
    @("273÷6=", "815÷3="),
    @("592÷3=", "164÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
